$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Cells.Item(33,8).Value = 418.03333
$ws.Cells.Item(33,9).Value = 305.65384
$ws.Cells.Item(33,11).Value = 305.65384
$ws.Cells.Item(33,13).Value = -76.65384
# row 40
$ws.Cells.Item(40,8).Value = 3565.125
$ws.Cells.Item(40,9).Value = 3733
$ws.Cells.Item(40,11).Value = 3733
$ws.Cells.Item(40,13).Value = -3558
# row 46
$ws.Cells.Item(46,8).Value = 201495.8
$ws.Cells.Item(46,10).Value = 502550
$ws.Cells.Item(46,12).Value = 1507650
$ws.Cells.Item(46,14).Value = -1507888
# row 60
$ws.Cells.Item(60,8).Value = 201495.8
$ws.Cells.Item(60,10).Value = 502550
$ws.Cells.Item(60,12).Value = 1507650
$ws.Cells.Item(60,14).Value = -1508618
# row 64
$ws.Cells.Item(64,8).Value = 55931.125
$ws.Cells.Item(64,9).Value = 71741.5
$ws.Cells.Item(64,11).Value = 71741.5
$ws.Cells.Item(64,13).Value = -71493.5
# row 67
$ws.Cells.Item(67,8).Value = 55931.125
$ws.Cells.Item(67,9).Value = 71741.5
$ws.Cells.Item(67,11).Value = 71741.5
$ws.Cells.Item(67,13).Value = -70883.5
# row 96
$ws.Cells.Item(96,8).Value = 876.6667
$ws.Cells.Item(96,9).Value = 646.6667
$ws.Cells.Item(96,11).Value = 1940.0001
$ws.Cells.Item(96,13).Value = -567.0001
# row 111
$ws.Cells.Item(111,8).Value = 5074.25
$ws.Cells.Item(111,9).Value = 3198.75
$ws.Cells.Item(111,11).Value = 9596.25
$ws.Cells.Item(111,13).Value = -6529.25
# row 132
$ws.Cells.Item(132,8).Value = 4004.4866
$ws.Cells.Item(132,10).Value = 8018.4
$ws.Cells.Item(132,12).Value = 24055.2
$ws.Cells.Item(132,14).Value = -29115.2
# row 135
$ws.Cells.Item(135,8).Value = 3150.4375
$ws.Cells.Item(135,9).Value = 3150.4375
$ws.Cells.Item(135,11).Value = 28353.9375
$ws.Cells.Item(135,13).Value = -25818.9375
# row 137
$ws.Cells.Item(137,8).Value = 6046.375
$ws.Cells.Item(137,9).Value = 6772.756
$ws.Cells.Item(137,10).Value = 1791.8572
$ws.Cells.Item(137,11).Value = 20318.268
$ws.Cells.Item(137,12).Value = 5375.571599999999
$ws.Cells.Item(137,13).Value = -17768.268
$ws.Cells.Item(137,14).Value = -10475.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Cells.Item(4,8).Value = 558.75
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(4,13).ClearContents()
# row 6
$ws.Cells.Item(6,8).Value = 10010000
$ws.Cells.Item(6,9).Value = 10010000
$ws.Cells.Item(6,11).Value = 10010000
$ws.Cells.Item(6,13).Value = -10009827
# row 9
$ws.Cells.Item(9,8).Value = 0
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(9,11).Value = 0
$ws.Cells.Item(9,13).ClearContents()
# row 20
$ws.Cells.Item(20,8).Value = 0
$ws.Cells.Item(20,9).Value = 0
$ws.Cells.Item(20,11).Value = 0
$ws.Cells.Item(20,13).ClearContents()
# row 32
$ws.Cells.Item(32,8).Value = 2964.0244
$ws.Cells.Item(32,9).Value = 3069.658
$ws.Cells.Item(32,11).Value = 3069.658
$ws.Cells.Item(32,13).Value = -2782.658
# row 44
$ws.Cells.Item(44,8).Value = 77999.2
$ws.Cells.Item(44,10).Value = 77999.2
$ws.Cells.Item(44,12).Value = 77999.2
$ws.Cells.Item(44,14).Value = -78975.2
# row 45
$ws.Cells.Item(45,8).Value = 4725.7085
$ws.Cells.Item(45,10).Value = 5743.222
$ws.Cells.Item(45,12).Value = 5743.222
$ws.Cells.Item(45,14).Value = -6497.222
# row 55
$ws.Cells.Item(55,8).Value = 45709.2
$ws.Cells.Item(55,9).Value = 9515.333000000001
$ws.Cells.Item(55,11).Value = 9515.333000000001
$ws.Cells.Item(55,13).Value = -9200.333000000001
# row 63
$ws.Cells.Item(63,8).Value = 6513.6665
$ws.Cells.Item(63,9).Value = 6513.6665
$ws.Cells.Item(63,10).Value = 0
$ws.Cells.Item(63,11).Value = 6513.6665
$ws.Cells.Item(63,12).Value = 0
$ws.Cells.Item(63,13).Value = -5827.6665
$ws.Cells.Item(63,14).ClearContents()
# row 66
$ws.Cells.Item(66,8).Value = 6513.6665
$ws.Cells.Item(66,9).Value = 6513.6665
$ws.Cells.Item(66,10).Value = 0
$ws.Cells.Item(66,11).Value = 32568.3325
$ws.Cells.Item(66,12).Value = 0
$ws.Cells.Item(66,13).Value = -29136.3325
$ws.Cells.Item(66,14).ClearContents()
# row 74
$ws.Cells.Item(74,8).Value = 2412.9666
$ws.Cells.Item(74,9).Value = 1305.56
$ws.Cells.Item(74,10).Value = 7950
$ws.Cells.Item(74,11).Value = 1305.56
$ws.Cells.Item(74,12).Value = 7950
$ws.Cells.Item(74,13).Value = -431.5599999999999
$ws.Cells.Item(74,14).Value = -9698
# row 77
$ws.Cells.Item(77,8).Value = 2412.9666
$ws.Cells.Item(77,9).Value = 1305.56
$ws.Cells.Item(77,10).Value = 7950
$ws.Cells.Item(77,11).Value = 6527.799999999999
$ws.Cells.Item(77,12).Value = 39750
$ws.Cells.Item(77,13).Value = -2159.799999999999
$ws.Cells.Item(77,14).Value = -48486
# row 80
$ws.Cells.Item(80,8).Value = 85000
$ws.Cells.Item(80,10).Value = 85000
$ws.Cells.Item(80,12).Value = 85000
$ws.Cells.Item(80,14).Value = -86996
# row 83
$ws.Cells.Item(83,8).Value = 85000
$ws.Cells.Item(83,10).Value = 85000
$ws.Cells.Item(83,12).Value = 255000
$ws.Cells.Item(83,14).Value = -264984
# row 132
$ws.Cells.Item(132,8).Value = 2306.7036
$ws.Cells.Item(132,9).Value = 1375.2609
$ws.Cells.Item(132,11).Value = 4125.7827
$ws.Cells.Item(132,13).Value = -1595.7827

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 140
$ws.Cells.Item(140,8).Value = 85709
$ws.Cells.Item(140,10).Value = 0
$ws.Cells.Item(140,12).Value = 0
$ws.Cells.Item(140,14).ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Cells.Item(7,8).Value = 13654.1875
$ws.Cells.Item(7,9).Value = 27068.375
$ws.Cells.Item(7,10).Value = 240
$ws.Cells.Item(7,11).Value = 27068.375
$ws.Cells.Item(7,12).Value = 240
$ws.Cells.Item(7,13).Value = -26955.375
$ws.Cells.Item(7,14).Value = -466
# row 139
$ws.Cells.Item(139,8).Value = 55498.25
$ws.Cells.Item(139,10).Value = 55498.25
$ws.Cells.Item(139,12).Value = 55498.25
$ws.Cells.Item(139,14).Value = -65778.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 51
$ws.Cells.Item(51,8).Value = 1929.4286
$ws.Cells.Item(51,10).Value = 2863
$ws.Cells.Item(51,12).Value = 8589
$ws.Cells.Item(51,14).Value = -9509
# row 63
$ws.Cells.Item(63,8).Value = 2250
$ws.Cells.Item(63,10).Value = 3000
$ws.Cells.Item(63,12).Value = 9000
$ws.Cells.Item(63,14).Value = -10498
# row 66
$ws.Cells.Item(66,8).Value = 2250
$ws.Cells.Item(66,10).Value = 3000
$ws.Cells.Item(66,12).Value = 27000
$ws.Cells.Item(66,14).Value = -34488
# row 132
$ws.Cells.Item(132,8).Value = 47920.453
$ws.Cells.Item(132,9).Value = 931.5
$ws.Cells.Item(132,10).Value = 74771.28999999999
$ws.Cells.Item(132,11).Value = 8383.5
$ws.Cells.Item(132,12).Value = 672941.61
$ws.Cells.Item(132,13).Value = -5853.5
$ws.Cells.Item(132,14).Value = -678001.61
# row 140
$ws.Cells.Item(140,8).Value = 8990.923000000001
$ws.Cells.Item(140,9).Value = 10829.685
$ws.Cells.Item(140,11).Value = 32489.055
$ws.Cells.Item(140,13).Value = -27309.055

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Cells.Item(122,8).Value = 19622.223
$ws.Cells.Item(122,9).Value = 25165.666
$ws.Cells.Item(122,10).Value = 16850.5
$ws.Cells.Item(122,11).Value = 75496.99800000001
$ws.Cells.Item(122,12).Value = 50551.5
$ws.Cells.Item(122,13).Value = -73046.99800000001
$ws.Cells.Item(122,14).Value = -55451.5
# row 132
$ws.Cells.Item(132,8).Value = 3161.7144
$ws.Cells.Item(132,9).Value = 2969.842
$ws.Cells.Item(132,11).Value = 8909.526
$ws.Cells.Item(132,13).Value = -6379.526

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 43
$ws.Cells.Item(43,8).Value = 16568.5
$ws.Cells.Item(43,9).Value = 16996.666
$ws.Cells.Item(43,10).Value = 13999.5
$ws.Cells.Item(43,11).Value = 16996.666
$ws.Cells.Item(43,12).Value = 13999.5
$ws.Cells.Item(43,13).Value = -16803.666
$ws.Cells.Item(43,14).Value = -14385.5
# row 55
$ws.Cells.Item(55,8).Value = 1688.5834
$ws.Cells.Item(55,9).Value = 344
$ws.Cells.Item(55,10).Value = 2649
$ws.Cells.Item(55,11).Value = 344
$ws.Cells.Item(55,12).Value = 2649
$ws.Cells.Item(55,13).Value = -171
$ws.Cells.Item(55,14).Value = -2995
# row 68
$ws.Cells.Item(68,8).Value = 4380.273
$ws.Cells.Item(68,10).Value = 4568.5
$ws.Cells.Item(68,12).Value = 4568.5
$ws.Cells.Item(68,14).Value = -6066.5
# row 71
$ws.Cells.Item(71,8).Value = 4380.273
$ws.Cells.Item(71,10).Value = 4568.5
$ws.Cells.Item(71,12).Value = 22842.5
$ws.Cells.Item(71,14).Value = -30330.5
# row 82
$ws.Cells.Item(82,8).Value = 2903.2
$ws.Cells.Item(82,10).Value = 1835.091
$ws.Cells.Item(82,12).Value = 1835.091
$ws.Cells.Item(82,14).Value = -2557.091
# row 85
$ws.Cells.Item(85,8).Value = 2903.2
$ws.Cells.Item(85,10).Value = 1835.091
$ws.Cells.Item(85,12).Value = 1835.091
$ws.Cells.Item(85,14).Value = -4331.091
# row 97
$ws.Cells.Item(97,8).Value = 19390.834
$ws.Cells.Item(97,10).Value = 19390.834
$ws.Cells.Item(97,12).Value = 19390.834
$ws.Cells.Item(97,14).Value = -21372.834
# row 132
$ws.Cells.Item(132,8).Value = 624118.8
$ws.Cells.Item(132,9).Value = 1147842.5
$ws.Cells.Item(132,11).Value = 3443527.5
$ws.Cells.Item(132,13).Value = -3440997.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 24
$ws.Cells.Item(24,8).Value = 26112.785
$ws.Cells.Item(24,10).Value = 26929.154
$ws.Cells.Item(24,12).Value = 26929.154
$ws.Cells.Item(24,14).Value = -27389.154
# row 40
$ws.Cells.Item(40,8).Value = 39663
$ws.Cells.Item(40,9).Value = 8999
$ws.Cells.Item(40,11).Value = 8999
$ws.Cells.Item(40,13).Value = -8850
# row 62
$ws.Cells.Item(62,8).Value = 273928.44
$ws.Cells.Item(62,9).Value = 950001
$ws.Cells.Item(62,10).Value = 3499.4
$ws.Cells.Item(62,11).Value = 950001
$ws.Cells.Item(62,12).Value = 3499.4
$ws.Cells.Item(62,13).Value = -949377
$ws.Cells.Item(62,14).Value = -4747.4
# row 65
$ws.Cells.Item(65,8).Value = 273928.44
$ws.Cells.Item(65,9).Value = 950001
$ws.Cells.Item(65,10).Value = 3499.4
$ws.Cells.Item(65,11).Value = 4750005
$ws.Cells.Item(65,12).Value = 17497
$ws.Cells.Item(65,13).Value = -4746885
$ws.Cells.Item(65,14).Value = -23737
# row 95
$ws.Cells.Item(95,8).Value = 86033830
$ws.Cells.Item(95,10).Value = 86033830
$ws.Cells.Item(95,12).Value = 86033830
$ws.Cells.Item(95,14).Value = -86039322
# row 132
$ws.Cells.Item(132,8).Value = 8798.016
$ws.Cells.Item(132,9).Value = 9264.383
$ws.Cells.Item(132,11).Value = 27793.149
$ws.Cells.Item(132,13).Value = -25263.149
